$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before row 21; existing rows 21-62 shift down to 22-63
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100107
$ws.Range("H21").Value = "Otros"
$ws.Range("I21").Value = 100107002
$ws.Range("J21").Value = "Chirimoya"
$ws.Range("K21").Value = "Cultivar IV Región"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 90
$ws.Range("N21").Value = 3000
$ws.Range("O21").Value = 3200
$ws.Range("P21").Value = 3111
$ws.Range("Q21").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R21").Value = "Provincia del Elquí"
$ws.Range("S21").Value = 3111
$ws.Range("T21").Value = 1
